$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the data range as Text so numeric-looking strings (e.g. "4.01")
# are preserved verbatim instead of being auto-converted to numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Apply the updated price / volume figures scraped on Fri Aug 16 04:47:59 UTC 2024.
$ws.Range("D2").Value = "58.162.51"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "2.574.16"
$ws.Range("E3").Value = "  -2.57%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "518.54"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "142.28"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "0.564"
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").Value = "2.588.59"
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("D10").Value = "6.78"
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").Value = "0.325"
$ws.Range("E12").Value = "  -3.94%  "
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").Value = "3.027.22"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").Value = "58.090.61"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").Value = "20.33"
$ws.Range("E16").Value = "  -2.43%  "
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").Value = "2.567.98"
$ws.Range("E18").Value = "  -2.98%  "
$ws.Range("D19").Value = "341.42"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").Value = "4.30"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("D21").Value = "10.26"
$ws.Range("E21").Value = "  -1.69%  "
$ws.Range("D22").Value = "6.34"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("D24").Value = "65.54"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("E25").Value = "  -1.47%  "
$ws.Range("E26").Value = "  -5.36%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").Value = "2.686.05"
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("D29").Value = "6.99"
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("D30").Value = "0.0₃0745"
$ws.Range("E30").Value = "  -6.20%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "6.24"
$ws.Range("E32").Value = "  -6.46%  "
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").Value = "18.69"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("D35").Value = "149.90"
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("D36").Value = "4.01"
$ws.Range("E36").Value = "  -3.01%  "
$ws.Range("E37").Value = "  -3.43%  "
$ws.Range("D38").Value = "0.865"
$ws.Range("E38").Value = "  -4.97%  "
$ws.Range("D39").Value = "35.95"
$ws.Range("E39").Value = "  -2.23%  "
$ws.Range("D40").Value = "0.835"
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("D41").Value = "1.45"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("D42").Value = "3.53"
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").Value = "269.81"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "10.67"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").Value = "0.0949"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("E47").Value = "  -3.09%  "
$ws.Range("E48").Value = "  -2.97%  "
$ws.Range("D49").Value = "0.0523"
$ws.Range("E49").Value = "  -2.72%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.973.14"
$ws.Range("E50").Value = "  -3.28%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "4.61"
$ws.Range("E51").Value = "  -0.80%  "

# Restore the default (unstyled) appearance now that the text values are locked in.
$dataRange.Style = "Normal"

